{"js": "// Split three single-run paragraphs (\"Programa\" PT/EN list and the\n// \"Bibliografia\" references) into multiple <w:t> runs joined by <w:br/>\n// line breaks, one per enumerated item / reference, without altering any\n// of the visible text itself.\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\n// Build the <w:r> markup for a run broken across several <w:t>/<w:br/>\n// segments. `segments` is an array of plain-text chunks; a <w:br/> is\n// placed between every pair of consecutive chunks. Runs keep whatever\n// <w:rPr> (e.g. italics) the original paragraph used.\nfunction buildBrokenRunOoxml(segments, rPrXml) {\n  const rPr = rPrXml || \"\";\n  const parts = segments.map((seg) => {\n    const preserve = /^\\s|\\s$/.test(seg) ? ' xml:space=\"preserve\"' : \"\";\n    return `<w:t${preserve}>${escapeXml(seg)}</w:t>`;\n  });\n  const body = parts.join(\"<w:br/>\");\n  return `<w:p><w:r>${rPr}${body}</w:r></w:p>`;\n}\n\nfunction wrapPackageOoxml(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    paragraphXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three target paragraphs by their distinctive text prefixes\n// rather than hard-coded indices, so the script is resilient to minor\n// structural differences.\nlet idxPtPrograma = -1;\nlet idxEnPrograma = -1;\nlet idxBibliografia = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"1) \u00d3ptica Geom\u00e9trica: conceitos b\u00e1sicos.\") === 0) {\n    idxPtPrograma = i;\n  } else if (t.indexOf(\"1) Geometrical Optics: basic concepts.\") === 0) {\n    idxEnPrograma = i;\n  } else if (t.indexOf(\"NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica.\") === 0) {\n    idxBibliografia = i;\n  }\n}\n\nif (idxPtPrograma === -1 || idxEnPrograma === -1 || idxBibliografia === -1) {\n  throw new Error(\"Could not locate one or more target paragraphs.\");\n}\n\n// 1) Portuguese \"Programa\" paragraph.\nconst ptSegments = [\n  \"1) \u00d3ptica Geom\u00e9trica: conceitos b\u00e1sicos. \",\n  \"2) Interfer\u00eancia: a experi\u00eancia de Young; coer\u00eancia; figuras de interfer\u00eancia; o interfer\u00f4metro de Michelson.\",\n  \"3) Difra\u00e7\u00e3o.\",\n  \"4) Polariza\u00e7\u00e3o.\",\n  \"5) Relatividade: os postulados da relatividade, as transforma\u00e7\u00f5es de Lorentz, simultaneidade, tempo e comprimento; momento linear, trabalho e energia;\",\n  \"6) Prim\u00f3rdios da teoria qu\u00e2ntica: a hip\u00f3tese de Plank; o efeito fotoel\u00e9trico, quantiza\u00e7\u00e3o do f\u00f3ton; ondas de De Broglie, o efeito Compton, a difra\u00e7\u00e3o de el\u00e9trons, interfer\u00eancia; \",\n  \"7) Princ\u00edpios b\u00e1sicos da mec\u00e2nica qu\u00e2ntica: o princ\u00edpio de incerteza; a equa\u00e7\u00e3o de Schr\u00f6dinger.\",\n];\n\n// 2) English \"Programa\" paragraph (italic run).\nconst enSegments = [\n  \"1) Geometrical Optics: basic concepts.\",\n  \"2) Interference: Young's experience; coherence; interference figures; the Michelson interferometer.\",\n  \"3) Diffraction.\",\n  \"4) Polarization.\",\n  \"5) Relativity: the postulates of relativity, Lorentz transformations, simultaneity, time and length; linear momentum, work and energy;\",\n  \"6) Early days of quantum theory: the hypothesis of Planck; the photoelectric effect, quantization of the photon; De Broglie waves, the Compton effect, the electron diffraction, interference;\",\n  \"7) Basic principles of quantum mechanics: the uncertainty principle; the Schr\u00f6dinger equation.\",\n];\n\n// 3) \"Bibliografia\" paragraph.\nconst bibSegments = [\n  \"NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica. Vol. 4, Edgard Blucher (2008).\",\n  \"RESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 4, LTC (2008).\",\n  \"TIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 4, LTC (2008).\",\n  \"SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica IV, Vol. 4, Pearson Addison Wesley (2009).\",\n  \"JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 4, Thomson Pioneira (2008).\",\n];\n\nconst ptOoxml = wrapPackageOoxml(buildBrokenRunOoxml(ptSegments, \"\"));\nconst enOoxml = wrapPackageOoxml(\n  buildBrokenRunOoxml(enSegments, \"<w:rPr><w:i/></w:rPr>\")\n);\nconst bibOoxml = wrapPackageOoxml(buildBrokenRunOoxml(bibSegments, \"\"));\n\nparagraphs.items[idxPtPrograma].getRange(\"Whole\").insertOoxml(ptOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\nparagraphs.items[idxEnPrograma].getRange(\"Whole\").insertOoxml(enOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\nparagraphs.items[idxBibliografia].getRange(\"Whole\").insertOoxml(bibOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split three single-run paragraphs (\"Programa\" PT/EN list and the\n# \"Bibliografia\" references) into multiple <w:t> runs joined by <w:br/>\n# line breaks, one per enumerated item / reference, without altering any\n# of the visible text itself.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$idxPtPrograma = 0\n$idxEnPrograma = 0\n$idxBibliografia = 0\n\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($idxPtPrograma -eq 0 -and $t.StartsWith(\"1) \u00d3ptica Geom\u00e9trica: conceitos b\u00e1sicos.\")) {\n        $idxPtPrograma = $i\n    } elseif ($idxEnPrograma -eq 0 -and $t.StartsWith(\"1) Geometrical Optics: basic concepts.\")) {\n        $idxEnPrograma = $i\n    } elseif ($idxBibliografia -eq 0 -and $t.StartsWith(\"NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica.\")) {\n        $idxBibliografia = $i\n    }\n}\n\nif ($idxPtPrograma -eq 0 -or $idxEnPrograma -eq 0 -or $idxBibliografia -eq 0) {\n    throw \"Could not locate one or more target paragraphs.\"\n}\n\n$ooxmlHeader = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$ooxmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# 1) Portuguese \"Programa\" paragraph.\n$ptOoxml = $ooxmlHeader + '<w:p><w:r><w:t xml:space=\"preserve\">1) \u00d3ptica Geom\u00e9trica: conceitos b\u00e1sicos. </w:t><w:br/><w:t>2) Interfer\u00eancia: a experi\u00eancia de Young; coer\u00eancia; figuras de interfer\u00eancia; o interfer\u00f4metro de Michelson.</w:t><w:br/><w:t>3) Difra\u00e7\u00e3o.</w:t><w:br/><w:t>4) Polariza\u00e7\u00e3o.</w:t><w:br/><w:t>5) Relatividade: os postulados da relatividade, as transforma\u00e7\u00f5es de Lorentz, simultaneidade, tempo e comprimento; momento linear, trabalho e energia;</w:t><w:br/><w:t xml:space=\"preserve\">6) Prim\u00f3rdios da teoria qu\u00e2ntica: a hip\u00f3tese de Plank; o efeito fotoel\u00e9trico, quantiza\u00e7\u00e3o do f\u00f3ton; ondas de De Broglie, o efeito Compton, a difra\u00e7\u00e3o de el\u00e9trons, interfer\u00eancia; </w:t><w:br/><w:t>7) Princ\u00edpios b\u00e1sicos da mec\u00e2nica qu\u00e2ntica: o princ\u00edpio de incerteza; a equa\u00e7\u00e3o de Schr\u00f6dinger.</w:t></w:r></w:p>' + $ooxmlFooter\n\n# 2) English \"Programa\" paragraph (italic run).\n$enOoxml = $ooxmlHeader + '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>1) Geometrical Optics: basic concepts.</w:t><w:br/><w:t>2) Interference: Young''s experience; coherence; interference figures; the Michelson interferometer.</w:t><w:br/><w:t>3) Diffraction.</w:t><w:br/><w:t>4) Polarization.</w:t><w:br/><w:t>5) Relativity: the postulates of relativity, Lorentz transformations, simultaneity, time and length; linear momentum, work and energy;</w:t><w:br/><w:t>6) Early days of quantum theory: the hypothesis of Planck; the photoelectric effect, quantization of the photon; De Broglie waves, the Compton effect, the electron diffraction, interference;</w:t><w:br/><w:t>7) Basic principles of quantum mechanics: the uncertainty principle; the Schr\u00f6dinger equation.</w:t></w:r></w:p>' + $ooxmlFooter\n\n# 3) \"Bibliografia\" paragraph.\n$bibOoxml = $ooxmlHeader + '<w:p><w:r><w:t>NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica. Vol. 4, Edgard Blucher (2008).</w:t><w:br/><w:t>RESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 4, LTC (2008).</w:t><w:br/><w:t>TIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 4, LTC (2008).</w:t><w:br/><w:t>SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica IV, Vol. 4, Pearson Addison Wesley (2009).</w:t><w:br/><w:t>JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 4, Thomson Pioneira (2008).</w:t></w:r></w:p>' + $ooxmlFooter\n\n$d.Paragraphs.Item($idxPtPrograma).Range.InsertXML($ptOoxml)\n$d.Paragraphs.Item($idxEnPrograma).Range.InsertXML($enOoxml)\n$d.Paragraphs.Item($idxBibliografia).Range.InsertXML($bibOoxml)\n"}
